$d = $word.ActiveDocument
$d.Tables(1).Cell(1, 1).Range.Text = '65÷7=9, 2'
$d.Tables(1).Cell(1, 2).Range.Text = '73÷2=36, 1'
$d.Tables(1).Cell(1, 3).Range.Text = '98÷8=12, 2'
$d.Tables(1).Cell(1, 4).Range.Text = '71÷6=11, 5'
$d.Tables(1).Cell(1, 5).Range.Text = '88÷7=12, 4'
$d.Tables(1).Cell(5, 1).Range.Text = '37÷8=4, 5'
$d.Tables(1).Cell(5, 2).Range.Text = '33÷9=3, 6'
$d.Tables(1).Cell(5, 3).Range.Text = '85÷9=9, 4'
$d.Tables(1).Cell(5, 4).Range.Text = '58÷5=11, 3'
$d.Tables(1).Cell(5, 5).Range.Text = '47÷7=6, 5'
$d.Tables(1).Cell(9, 1).Range.Text = '90÷9=10, 0'
$d.Tables(1).Cell(9, 2).Range.Text = '29÷8=3, 5'
$d.Tables(1).Cell(9, 3).Range.Text = '15÷8=1, 7'
$d.Tables(1).Cell(9, 4).Range.Text = '40÷6=6, 4'
$d.Tables(1).Cell(9, 5).Range.Text = '78÷5=15, 3'
$d.Tables(1).Cell(13, 1).Range.Text = '69÷8=8, 5'
$d.Tables(1).Cell(13, 2).Range.Text = '53÷7=7, 4'
$d.Tables(1).Cell(13, 3).Range.Text = '32÷5=6, 2'
$d.Tables(1).Cell(13, 4).Range.Text = '47÷2=23, 1'
$d.Tables(1).Cell(13, 5).Range.Text = '18÷6=3, 0'
$d.Tables(1).Cell(17, 1).Range.Text = '41÷5=8, 1'
$d.Tables(1).Cell(17, 2).Range.Text = '61÷9=6, 7'
$d.Tables(1).Cell(17, 3).Range.Text = '85÷5=17, 0'
$d.Tables(1).Cell(17, 4).Range.Text = '43÷4=10, 3'
$d.Tables(1).Cell(17, 5).Range.Text = '17÷3=5, 2'
